# Apply updated values to subject32.xlsx connectivity matrix (GroupName2)
# Each assignment sets a specific cell's value to the new value from the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = 0.83367850075631211
$ws.Range("P1").Value = 0.80546468036762431
$ws.Range("C2").Value = 0.79594628773323817
$ws.Range("H2").Value = 0.97421105818463105
$ws.Range("AR3").Value = 0.73726103884524596
$ws.Range("AW3").Value = 0.91725930704801106
$ws.Range("O5").Value = 0.73385125611160718
$ws.Range("AX5").Value = 0.97758810178651356
$ws.Range("V7").Value = 0.81990362286020391
$ws.Range("J8").Value = 0.85508501224711653
$ws.Range("P8").Value = 0.9049259404877239
$ws.Range("AN8").Value = 0.61342890890388668
$ws.Range("AO9").Value = 0.85536428381349483
$ws.Range("AR9").Value = 0.99115334326206295
$ws.Range("BB10").Value = 0.97771854672759684
$ws.Range("N11").Value = 0.86240590845704501
$ws.Range("AB11").Value = 0.57306516603682645
$ws.Range("AL11").Value = 0.93803239341811362
$ws.Range("BM11").Value = 0.84737778675366782
$ws.Range("F12").Value = 0.76634458077853918
$ws.Range("H12").Value = 0.90763049369300286
$ws.Range("Z12").Value = 0.71055214763513308
$ws.Range("N13").Value = 0.57419693411036465
$ws.Range("BD13").Value = 0.79423946609522411
$ws.Range("AH14").Value = 0.76721064781839754
$ws.Range("AD15").Value = 0.6060931434771919
$ws.Range("AE15").Value = 0.86015590200119929
$ws.Range("R16").Value = 0.65736593158940959
$ws.Range("AD16").Value = 0.99096292763283245
$ws.Range("D17").Value = 0.92116389330060544
$ws.Range("BD17").Value = 0.92915662703470137
$ws.Range("AG18").Value = 0.96910095400005014
$ws.Range("AM18").Value = 0.86171043722912644
$ws.Range("F19").Value = 0.84378895461969061
$ws.Range("I19").Value = 0.93730289642844378
$ws.Range("AD19").Value = 0.58844452144416981
$ws.Range("AG19").Value = 0.81902583692844089
$ws.Range("U20").Value = 0.6737704604755046
$ws.Range("Y20").Value = 0.73767985289979721
$ws.Range("Z20").Value = 0.96534313263339122
$ws.Range("AB20").Value = 0.99412515001929602
$ws.Range("V21").Value = 0.90073244317818135
$ws.Range("AM21").Value = 0.90415829313411566
$ws.Range("AO21").Value = 0.86562275829207636
$ws.Range("BH21").Value = 0.76150889410070077
$ws.Range("BF22").Value = 0.60242592045525245
$ws.Range("O24").Value = 0.83503087957812172
$ws.Range("G25").Value = 0.83171515488565095
$ws.Range("W25").Value = 0.86079322054775731
$ws.Range("A27").Value = 0.63436833041432206
$ws.Range("W27").Value = 0.94874398329184284
$ws.Range("AH27").Value = 0.7130945816562646
$ws.Range("BI27").Value = 0.58303607835360516
$ws.Range("BK27").Value = 0.69232688779718221
$ws.Range("L28").Value = 0.70622859463837351
$ws.Range("AC28").Value = 0.96869324992318684
$ws.Range("A29").Value = 0.93949500863609425
$ws.Range("D29").Value = 0.63223111923003139
$ws.Range("AJ29").Value = 0.86403905607353781
$ws.Range("AA30").Value = 0.5626727159380156
$ws.Range("AU30").Value = 0.82875188333669492
$ws.Range("AM33").Value = 0.75679702933024195
$ws.Range("BK33").Value = 0.80732535837513564
$ws.Range("B34").Value = 0.96143626262092519
$ws.Range("Y34").Value = 0.79027900994043221
$ws.Range("F35").Value = 0.88784415334368882
$ws.Range("AJ35").Value = 0.62974611716068374
$ws.Range("AV35").Value = 0.98826831017844574
$ws.Range("BJ35").Value = 0.96322965188121157
$ws.Range("H36").Value = 0.89068434726611723
$ws.Range("AH36").Value = 0.88402039759753259
$ws.Range("B37").Value = 0.92082721716303539
$ws.Range("D38").Value = 0.73354284620833221
$ws.Range("AO38").Value = 0.84642041730904105
$ws.Range("BE38").Value = 0.84707039320544597
$ws.Range("AL39").Value = 0.8969379959711532
$ws.Range("O40").Value = 0.99093298848891109
$ws.Range("C41").Value = 0.97626921097861763
$ws.Range("BB41").Value = 0.63805913676292292
$ws.Range("X42").Value = 0.94817080582581115
$ws.Range("AZ42").Value = 0.9989875595329506
$ws.Range("BL42").Value = 0.76860999130379448
$ws.Range("BM43").Value = 0.6046001550914617
$ws.Range("AT44").Value = 0.80813794387739368
$ws.Range("AA45").Value = 0.82772193846915099
$ws.Range("AF45").Value = 0.71270109438241902
$ws.Range("D46").Value = 0.98032062064436665
$ws.Range("AK46").Value = 0.80863524987317703
$ws.Range("AM46").Value = 0.76999348838405179
$ws.Range("AS46").Value = 0.83364736338983958
$ws.Range("A47").Value = 0.87638582686951605
$ws.Range("G47").Value = 0.84626384899774409
$ws.Range("O47").Value = 0.90340968562550406
$ws.Range("Z47").Value = 0.91903416787581049
$ws.Range("AP47").Value = 0.86161553924925061
$ws.Range("AQ47").Value = 0.97498812610696217
$ws.Range("AX49").Value = 0.88166014967242878
$ws.Range("BA49").Value = 0.81015613351958304
$ws.Range("AY50").Value = 0.95064327884037358
$ws.Range("X51").Value = 0.97299583141173906
$ws.Range("BK52").Value = 0.85107544866186391
$ws.Range("M53").Value = 0.98761723503726151
$ws.Range("S53").Value = 0.73323673561573699
$ws.Range("AE53").Value = 0.58068715132733673
$ws.Range("BC54").Value = 0.98570991322687418
$ws.Range("Q55").Value = 0.81597068187300525
$ws.Range("R55").Value = 0.83118288868164003
$ws.Range("AQ55").Value = 0.88156564213416488
$ws.Range("AY55").Value = 0.63382531321005298
$ws.Range("BI55").Value = 0.9755513849406996
$ws.Range("D56").Value = 0.91636570012272733
$ws.Range("R57").Value = 0.97074011468225418
$ws.Range("AO57").Value = 0.87901791050135181
$ws.Range("BD57").Value = 0.95112478785987387
$ws.Range("AF58").Value = 0.80209297741038488
$ws.Range("AP58").Value = 0.83887678565096313
$ws.Range("AT59").Value = 0.88827270391454249
$ws.Range("K60").Value = 0.61055255787522578
$ws.Range("AV60").Value = 0.7546068631805295
$ws.Range("BF60").Value = 0.97431467870386301
$ws.Range("M61").Value = 0.73863536381774819
$ws.Range("BG61").Value = 0.85040411114324255
$ws.Range("P62").Value = 0.9588052620894304
$ws.Range("S62").Value = 0.99413243016330155
$ws.Range("AO62").Value = 0.90616360763554415
$ws.Range("AY62").Value = 0.95323794512470084
$ws.Range("AN64").Value = 0.92040165894850567
$ws.Range("BM64").Value = 0.89615175216663168
$ws.Range("J65").Value = 0.78282305644365979
$ws.Range("BN65").Value = 0.7741895475273517
$ws.Range("Y66").Value = 0.97766456123837209
$ws.Range("AA66").Value = 0.65706045045956807
$ws.Range("AE67").Value = 0.73501209914911303
$ws.Range("AX67").Value = 0.94745349836883186
$ws.Range("AY68").Value = 0.99311489113998086
$ws.Range("BH68").Value = 0.93789911440514684
